$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45310 -> 45311, i.e. 2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# Update prices in the "Para CARGA" price list (column D)
$ws.Range("D23").Value = 693
$ws.Range("D24").Value = 785
$ws.Range("D25").Value = 954
$ws.Range("D26").Value = 1066
$ws.Range("D27").Value = 1196

# Update prices in the "Para DESCARGA" price list (column D)
$ws.Range("D34").Value = 577
$ws.Range("D35").Value = 725.1
